$wb = $excel.ActiveWorkbook

# --- SolverSettings sheet: append a new "include_RPS" settings row ---
$ws = $wb.Worksheets.Item("SolverSettings")

$ws.Range("A10").Value = "include_RPS"
$ws.Range("B10:G10").Value = "N"

# Make SolverSettings the active sheet/tab and select H10 on it
# (mirrors the new selection + tabSelected state from the authored edit).
$ws.Activate()
$ws.Range("H10").Select()

# --- Connections sheet: it is no longer the active/selected tab ---
# (Activating SolverSettings above already moves "tabSelected" away from
# Connections; nothing else on that sheet changes.)
